$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AfterShip")

$data = @(
    @("TEST_GDPR", "dbschenker-se", "Delivered", "12345", "", "2026-02-08T16:46:34+04:30", ",Afghanistan", "2026-02-08T12:25:31+00:00", "SHIPMENT_TITLE", '{"OrderID": "12345"}'),
    @("TEST_TRACKING_DECEMBER", "dhl", "Delivered", "", "", "2025-12-31T23:59:59", "", "2026-02-08T11:24:55+00:00", "TEST_Tracking_December", "{}"),
    @("TEST_TRACKING", "kn", "Delivered", "", "", "2026-02-07T23:59:59", "", "2026-02-08T11:24:18+00:00", "ZFRE", "{}"),
    @("ITD-0-12345678", "testing-courier", "Delivered", "", "", "2026-02-08T07:07:24-04:00", "1000 W Test Ln, Austin, Texas, 11111", "2026-02-08T11:07:24+00:00", "ITD-0-12345678", "{}")
)

# Columns that must remain text even when they look numeric
$textColumns = @(4)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        if ($textColumns -contains $colIndex -and $val -ne "") {
            # Prefix with an apostrophe so numeric-looking text (e.g. order ids)
            # is stored as text, same as typing '12345 in Excel.
            $cell.Value = "'" + $val
        }
        else {
            $cell.Value = $val
        }
        $colIndex++
    }
    $rowIndex++
}
